$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.939.21'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.00%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.905.57'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.31%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7968'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.21%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.26%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3117'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.89%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.44'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.23%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07002'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.99%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07995'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.12%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.904.46'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.30%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7381'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.39%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.169'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.74%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.47'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.06%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.954.43'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.864'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.08%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.22%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007766'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.51%  '

$ws.Range('E21').Value = '  +0.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.151.99'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.906'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.36%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '167.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.97%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.206'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.76%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1435'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.72%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.29%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.060'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.96%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.354'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.76%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.513'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.36%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.284'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.16%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05575'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.059'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.65%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.268'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.63%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7316'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.78%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.719'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01925'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.37%  '

$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4402'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.51%  '

$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.04'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.70%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.983'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.47%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.10%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8369'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.27%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.884'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.05%  '

$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.64'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.34%  '

$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.553'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.47%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.705'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.78%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '978.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.68%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.057.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.50%  '

$ws.Range('E51').Value = '  -0.58%  '
